# CIERRE 23 MAYO 22
# Apply updates to the "REMISIONES MAYO 2022" sheet (5th tab): fill in
# payment dates / amounts (columns F, G) and some missing remision rows
# (A, D, E) for the credit tracking table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Row 6-9: existing rows just get their "Fecha de pago" (F) and
#     "IMPORTE D/PAGO" (G) filled in (fully paid -> H recomputes to 0) ---
$ws.Cells.Item(6, 6).Value = 44695
$ws.Cells.Item(6, 7).Value = 54116

$ws.Cells.Item(7, 6).Value = 44688
$ws.Cells.Item(7, 7).Value = 6018

$ws.Cells.Item(8, 6).Value = 44687
$ws.Cells.Item(8, 7).Value = 15797

$ws.Cells.Item(9, 6).Value = 44687
$ws.Cells.Item(9, 7).Value = 1362

# --- Row 10: new remision entry (HERRADURA GUSTAVO) ---
$ws.Cells.Item(10, 1).Value = 44687
$ws.Cells.Item(10, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(10, 5).Value = 19666
$ws.Cells.Item(10, 6).Value = 44689
$ws.Cells.Item(10, 7).Value = 19666

# --- Row 11 ---
$ws.Cells.Item(11, 1).Value = 44688
$ws.Cells.Item(11, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(11, 5).Value = 24471
$ws.Cells.Item(11, 6).Value = 44688
$ws.Cells.Item(11, 7).Value = 24471

# --- Row 12 ---
$ws.Cells.Item(12, 1).Value = 44688
$ws.Cells.Item(12, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(12, 5).Value = 10200
$ws.Cells.Item(12, 6).Value = 44689
$ws.Cells.Item(12, 7).Value = 10200

# --- Row 13 ---
$ws.Cells.Item(13, 1).Value = 44689
$ws.Cells.Item(13, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(13, 5).Value = 15259
$ws.Cells.Item(13, 6).Value = 44691
$ws.Cells.Item(13, 7).Value = 15259

# --- Row 14 ---
$ws.Cells.Item(14, 1).Value = 44689
$ws.Cells.Item(14, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(14, 5).Value = 10231
$ws.Cells.Item(14, 6).Value = 44691
$ws.Cells.Item(14, 7).Value = 10231

# --- Row 15: MICH, not yet paid (F/G stay blank, H keeps the full saldo) ---
$ws.Cells.Item(15, 1).Value = 44691
$ws.Cells.Item(15, 4).Value = "MICH"
$ws.Cells.Item(15, 5).Value = 516

# --- Row 16 ---
$ws.Cells.Item(16, 1).Value = 44691
$ws.Cells.Item(16, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(16, 5).Value = 10876
$ws.Cells.Item(16, 6).Value = 44694
$ws.Cells.Item(16, 7).Value = 10876

# --- Row 17 ---
$ws.Cells.Item(17, 1).Value = 44691
$ws.Cells.Item(17, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(17, 5).Value = 10282
$ws.Cells.Item(17, 6).Value = 44694
$ws.Cells.Item(17, 7).Value = 10282

# --- Row 18: OBRADOR, not yet paid ---
$ws.Cells.Item(18, 1).Value = 44693
$ws.Cells.Item(18, 4).Value = "OBRADOR"
$ws.Cells.Item(18, 5).Value = 464

# --- Row 19 ---
$ws.Cells.Item(19, 1).Value = 44694
$ws.Cells.Item(19, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(19, 5).Value = 23516
$ws.Cells.Item(19, 6).Value = 44698
$ws.Cells.Item(19, 7).Value = 23516

# --- Row 20 ---
$ws.Cells.Item(20, 1).Value = 44694
$ws.Cells.Item(20, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(20, 5).Value = 8894
$ws.Cells.Item(20, 6).Value = 44695
$ws.Cells.Item(20, 7).Value = 8894

# --- Row 21 ---
$ws.Cells.Item(21, 1).Value = 44695
$ws.Cells.Item(21, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(21, 5).Value = 10200
$ws.Cells.Item(21, 6).Value = 44696
$ws.Cells.Item(21, 7).Value = 10200

# --- Row 22 ---
$ws.Cells.Item(22, 1).Value = 44696
$ws.Cells.Item(22, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(22, 5).Value = 11179
$ws.Cells.Item(22, 6).Value = 44698
$ws.Cells.Item(22, 7).Value = 11179

# --- Row 23: not yet paid ---
$ws.Cells.Item(23, 1).Value = 44698
$ws.Cells.Item(23, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(23, 5).Value = 11923

# --- Row 24: OBRADOR, not yet paid ---
$ws.Cells.Item(24, 1).Value = 44698
$ws.Cells.Item(24, 4).Value = "OBRADOR"
$ws.Cells.Item(24, 5).Value = 410

# --- Row 25 ---
$ws.Cells.Item(25, 1).Value = 44698
$ws.Cells.Item(25, 4).Value = "HERRADURA GUSTAVO"
$ws.Cells.Item(25, 5).Value = 6436
$ws.Cells.Item(25, 6).Value = 44699
$ws.Cells.Item(25, 7).Value = 6436

# --- Row 26: OBRADOR, not yet paid ---
$ws.Cells.Item(26, 1).Value = 44699
$ws.Cells.Item(26, 4).Value = "OBRADOR"
$ws.Cells.Item(26, 5).Value = 4445

# Scroll/selection state as left by the author at closing time
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E27").Select()
